$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 159, shifting existing rows 159-188 down to 160-189
$ws.Rows("159").Insert()

# Fill in the new row 159 with the new data record
$ws.Cells.Item(159, 1).Value = 3
$ws.Cells.Item(159, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(159, 3).Value = 'Coquimbo'
$ws.Cells.Item(159, 4).Value = 44694
$ws.Cells.Item(159, 5).Value = 5
$ws.Cells.Item(159, 6).Value = 100112030
$ws.Cells.Item(159, 7).Value = 'Poroto granado'
$ws.Cells.Item(159, 8).Value = 'Sin especificar'
$ws.Cells.Item(159, 9).Value = 'Primera'
$ws.Cells.Item(159, 10).Value = 73
$ws.Cells.Item(159, 11).Value = 22500
$ws.Cells.Item(159, 12).Value = 23000
$ws.Cells.Item(159, 13).Value = 22760
$ws.Cells.Item(159, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(159, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(159, 16).Value = 910
$ws.Cells.Item(159, 17).Value = 25
$ws.Cells.Item(159, 18).Value = 'Hortaliza'
